$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "266.85"
    4  = "6.344"
    5  = "0.06179"
    6  = "3.600"
    7  = "6.678"
    8  = "1.391"
    9  = "0.8298"
    10 = "0.01363"
    11 = "0.1609"
    12 = "0.08220"
    13 = "0.03399"
    14 = "0.03157"
    15 = "0.09270"
    16 = "3.923"
    17 = "0.001704"
    18 = "0.04848"
    19 = "0.006308"
    20 = "0.005379"
    21 = "0.001094"
    23 = "3.763"
    24 = "2.369"
    26 = "0.1246"
    27 = "0.0002692"
    40 = "0.04666"
    41 = "0.006918"
    42 = "0.1154"
    43 = "0.003473"
    44 = "0.01228"
    45 = "0.00006288"
    46 = "0.00000000753"
    47 = "0.7916"
    48 = "0.1569"
    49 = "0.00002107"
    50 = "0.01244"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.Style = "Normal"
}
